$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 78
$ws.Cells.Item(78, 1).Value = 45939
$ws.Cells.Item(78, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(78, 3).Value = 886.94400000000007
$ws.Cells.Item(78, 4).Value = 1280.3550000000002
$ws.Cells.Item(78, 5).Value = 446.76000000000005
$ws.Cells.Item(78, 6).Value = 619.61300000000006
$ws.Cells.Item(78, 7).Value = 540.86599999999999
$ws.Cells.Item(78, 8).Value = 688.09899999999993
$ws.Cells.Item(78, 9).Value = 585.82299999999998
$ws.Cells.Item(78, 10).Value = 335.48500000000001
$ws.Cells.Item(78, 11).Value = 134.5
$ws.Cells.Item(78, 12).Value = 436.24400000000003
$ws.Cells.Item(78, 13).Value = 158.55500000000001
$ws.Cells.Item(78, 14).Value = 263.22000000000003
$ws.Cells.Item(78, 15).Value = 997.20699999999999
$ws.Cells.Item(78, 16).Value = 1514.5870000000002
$ws.Cells.Item(78, 17).Value = 555.029
$ws.Cells.Item(78, 18).Value = 567.91700000000014
$ws.Cells.Item(78, 19).Value = 332.71099999999996
$ws.Cells.Item(78, 20).Value = 94.455999999999989
$ws.Cells.Item(78, 21).Value = 145.66
$ws.Cells.Item(78, 22).Value = 171.71999999999997
$ws.Cells.Item(78, 23).Value = 69.13
$ws.Cells.Item(78, 24).Value = 112.30000000000001
$ws.Cells.Item(78, 25).Value = 21.3
$ws.Cells.Item(78, 26).Value = 28.75

# Row 79
$ws.Cells.Item(79, 1).Value = 45939
$ws.Cells.Item(79, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(79, 3).Value = 681.17899999999986
$ws.Cells.Item(79, 4).Value = 526.1640000000001
$ws.Cells.Item(79, 5).Value = 31.5
$ws.Cells.Item(79, 6).Value = 10.188000000000001
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 129.54599999999999
$ws.Cells.Item(79, 9).Value = 196.81299999999999
$ws.Cells.Item(79, 10).Value = 344.14500000000004
$ws.Cells.Item(79, 11).Value = 327.57900000000001
$ws.Cells.Item(79, 12).Value = 132.387
$ws.Cells.Item(79, 13).Value = 244.458
$ws.Cells.Item(79, 14).Value = 211.197
$ws.Cells.Item(79, 15).Value = 486.00099999999998
$ws.Cells.Item(79, 16).Value = 645.1049999999999
$ws.Cells.Item(79, 17).Value = 238.88500000000002
$ws.Cells.Item(79, 18).Value = 369.64000000000004
$ws.Cells.Item(79, 19).Value = 148.09
$ws.Cells.Item(79, 20).Value = 195.75
$ws.Cells.Item(79, 21).Value = 68.329000000000008
$ws.Cells.Item(79, 22).Value = 25.774000000000001
$ws.Cells.Item(79, 23).Value = 0
$ws.Cells.Item(79, 24).Value = 29.094000000000001
$ws.Cells.Item(79, 25).Value = 0
$ws.Cells.Item(79, 26).Value = 99.465000000000003

# Row 80
$ws.Cells.Item(80, 1).Value = 45940
$ws.Cells.Item(80, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(80, 3).Value = 769.63999999999987
$ws.Cells.Item(80, 4).Value = 1456.8319999999999
$ws.Cells.Item(80, 5).Value = 245.96500000000003
$ws.Cells.Item(80, 6).Value = 335.91700000000003
$ws.Cells.Item(80, 7).Value = 626.51999999999987
$ws.Cells.Item(80, 8).Value = 591.226
$ws.Cells.Item(80, 9).Value = 517.13599999999997
$ws.Cells.Item(80, 10).Value = 265.03699999999998
$ws.Cells.Item(80, 11).Value = 180.68700000000001
$ws.Cells.Item(80, 12).Value = 279.43400000000003
$ws.Cells.Item(80, 13).Value = 144.38299999999998
$ws.Cells.Item(80, 14).Value = 172.60999999999999
$ws.Cells.Item(80, 15).Value = 921.97100000000012
$ws.Cells.Item(80, 16).Value = 1532.7189999999998
$ws.Cells.Item(80, 17).Value = 586.82000000000005
$ws.Cells.Item(80, 18).Value = 437.27399999999994
$ws.Cells.Item(80, 19).Value = 555.17800000000011
$ws.Cells.Item(80, 20).Value = 127.57000000000001
$ws.Cells.Item(80, 21).Value = 61.260000000000005
$ws.Cells.Item(80, 22).Value = 232.98
$ws.Cells.Item(80, 23).Value = 24.68
$ws.Cells.Item(80, 24).Value = 106.559
$ws.Cells.Item(80, 25).Value = 50.89
$ws.Cells.Item(80, 26).Value = 33.5

# Row 81
$ws.Cells.Item(81, 1).Value = 45940
$ws.Cells.Item(81, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(81, 3).Value = 407.14000000000004
$ws.Cells.Item(81, 4).Value = 502.07700000000006
$ws.Cells.Item(81, 5).Value = 36.07
$ws.Cells.Item(81, 6).Value = 189.82400000000001
$ws.Cells.Item(81, 7).Value = 68.289000000000001
$ws.Cells.Item(81, 8).Value = 159.73500000000001
$ws.Cells.Item(81, 9).Value = 150.226
$ws.Cells.Item(81, 10).Value = 249.48699999999994
$ws.Cells.Item(81, 11).Value = 351.99800000000005
$ws.Cells.Item(81, 12).Value = 250.45600000000005
$ws.Cells.Item(81, 13).Value = 182.54599999999999
$ws.Cells.Item(81, 14).Value = 298.392
$ws.Cells.Item(81, 15).Value = 267.12300000000005
$ws.Cells.Item(81, 16).Value = 729.20800000000008
$ws.Cells.Item(81, 17).Value = 311.584
$ws.Cells.Item(81, 18).Value = 223.46100000000001
$ws.Cells.Item(81, 19).Value = 159.32400000000001
$ws.Cells.Item(81, 20).Value = 257.14
$ws.Cells.Item(81, 21).Value = 81.396000000000001
$ws.Cells.Item(81, 22).Value = 10.132999999999999
$ws.Cells.Item(81, 23).Value = 61.464999999999989
$ws.Cells.Item(81, 24).Value = 106.045
$ws.Cells.Item(81, 25).Value = 12.273
$ws.Cells.Item(81, 26).Value = 89.150999999999996

$ws.Range("F84").Select() | Out-Null

Write-Output "done"